# Update gh-pages output data (ticket/view counts refreshed at 456a3b4)
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 3806
$ws1.Range("F7").Value  = 188
$ws1.Range("G13").Value = 45
$ws1.Range("F21").Value = 3345
$ws1.Range("F22").Value = 5683
$ws1.Range("F28").Value = 3325
$ws1.Range("F33").Value = 514
$ws1.Range("F38").Value = 109
$ws1.Range("F43").Value = 30

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 91

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 3806
$ws4.Range("F7").Value  = 188
$ws4.Range("F10").Value = 91
$ws4.Range("G14").Value = 45
$ws4.Range("F22").Value = 3345
$ws4.Range("F23").Value = 5683
$ws4.Range("F29").Value = 3325
$ws4.Range("F34").Value = 514
$ws4.Range("F39").Value = 109
$ws4.Range("F44").Value = 30
